## Fruta / hortaliza, semanal
## Insert this week's (2023-12-07 / serial 45267) price rows for
## "Terminal La Palmera de La Serena" - Chirimoya, at the top of the
## data block (new rows 378-381), shifting the rest of the data down
## by 4 rows. Also adds a new "Tercera" quality tier for this week.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room: insert 4 new blank rows starting at row 378 (old row 378
# and everything below shifts down by 4, to 382+).
$ws.Rows("378:381").Insert()

# Helper data, one row per new record (maps to final rows 378-381).
$newRows = @(
    @{ Row=378; Qty=360; Min=17000; Max=18000; Avg=17500; Kg=1750;
       Quality="Especial" },
    @{ Row=379; Qty=500; Min=14000; Max=15000; Avg=14500; Kg=1450;
       Quality="Primera" },
    @{ Row=380; Qty=400; Min=10000; Max=11000; Avg=10500; Kg=1050;
       Quality="Segunda" },
    @{ Row=381; Qty=240; Min=7000;  Max=8000;  Avg=7500;  Kg=750;
       Quality="Tercera" }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value  = 8
    $ws.Cells.Item($row, 2).Value  = "Terminal La Palmera de La Serena"
    $ws.Cells.Item($row, 3).Value  = "Coquimbo"
    $ws.Cells.Item($row, 4).Value  = "2023-12-07"
    $ws.Cells.Item($row, 5).Value  = 4
    $ws.Cells.Item($row, 6).Value  = "Fruta"
    $ws.Cells.Item($row, 7).Value  = 100107
    $ws.Cells.Item($row, 8).Value  = "Otros"
    $ws.Cells.Item($row, 9).Value  = 100107002
    $ws.Cells.Item($row, 10).Value = "Chirimoya"
    $ws.Cells.Item($row, 11).Value = "Cultivar IV Región"
    $ws.Cells.Item($row, 12).Value = $r.Quality
    $ws.Cells.Item($row, 13).Value = $r.Qty
    $ws.Cells.Item($row, 14).Value = $r.Min
    $ws.Cells.Item($row, 15).Value = $r.Max
    $ws.Cells.Item($row, 16).Value = $r.Avg
    $ws.Cells.Item($row, 17).Value = "`$/bandeja 10 kilos"
    $ws.Cells.Item($row, 18).Value = "Provincia de Limarí"
    $ws.Cells.Item($row, 19).Value = $r.Kg
    $ws.Cells.Item($row, 20).Value = 10
}
